$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Data edits -------------------------------------------------------
# Update the two driving input cells; every other changed cell on the
# sheet (C32, D32, E32, D35:E98, and all of the shared-formula ranges)
# is a formula that depends -- directly or transitively, through
# COMBIN($C$32, ...) and CEILING.MATH/FLOOR.MATH -- on these two inputs,
# so it recalculates automatically once they change.
$ws.Range("D30").Value = 1
$ws.Range("B32").Value = 21

# --- View-state edits ---------------------------------------------------
# Re-create the saved scroll position / active selection.
$ws.Activate()
$excel.Goto($ws.Range("A25"), $true)
$ws.Range("H36").Select()
